# Insert a new player row ("Elfman Strauss") into the player_stats sheet.
# The new row is inserted right after the "Jellal Fernandes" row (row 54),
# i.e. becomes the new row 55, pushing every subsequent row down by one.
#
# Because a native Range/Rows Insert() on this engine creates brand-new
# style (cellXfs) entries (and drops the custom row height), we instead
# shift the data for rows 55-100 down to 56-101 manually cell-by-cell
# (which re-uses the pre-existing style already applied to those cells),
# then populate row 55 with the new player's data, and finally fix up the
# style + row height for the one truly new row (101).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J")

# Row 101 does not exist yet; stamp it with the same formatting used by
# every other data row (taken from the last existing row, 100) before we
# start moving data around.
$ws.Range("A100:J100").Copy()
$ws.Range("A101:J101").PasteSpecial(-4122)

# Shift rows 100 down to 55 upward by one position (process bottom-up so
# we never clobber a source row before it has been read). The custom row
# height travels together with the row's data.
for ($r = 100; $r -ge 55; $r--) {
    $dest = $r + 1
    foreach ($col in $cols) {
        $ws.Range("$col$dest").Value = $ws.Range("$col$r").Value2
    }
    $ws.Rows.Item($dest).RowHeight = $ws.Rows.Item($r).RowHeight
}

# Populate the newly freed row 55 with Elfman Strauss' data.
$ws.Range("A55").Value = "Elfman Strauss"
$ws.Range("B55").Value = 0
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 41
$ws.Range("H55").Value = "Midfielder"
$ws.Range("I55").Value = "ft"
$ws.Range("J55").Value = "ft_5"

Write-Host "Row inserted"
